$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (legmaxROM indices) - columns B:E now represent reps at 15/16 deg
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - updated meanEMG values; C2/E2 have no data at these angles
$ws.Range("B2").Value = 2.3772960904422913
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.479341849987208
$ws.Range("E2").ClearContents()

# Row 3 (STR) - updated meanEMG values
$ws.Range("B3").Value = 2.2032100745536449
$ws.Range("C3").Value = -0.75226107008933984
$ws.Range("D3").Value = 3.6460566198073323
$ws.Range("E3").Value = -0.096784206100009193

# Update the saved selection to the updated data range
$ws.Range("B1:E3").Select()

$wb.Save()
